$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the userID values for existing rows (upload was assigning wrong IDs)
$ws.Range("A2").Value = 7
$ws.Range("A3").Value = 8

# rateType (H) for row 3 corrected; keep it text like the rest of that column
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "23"
$ws.Range("H3").Style = "Normal"

# Append the new user row that failed to upload previously
$ws.Range("A4").Value = 9
$ws.Range("B4").Value = "joe"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "42132"
$ws.Range("C4").Style = "Normal"

$ws.Range("D4").Value = "Graham Gibson"
$ws.Range("E4").Value = "Queen's University"
$ws.Range("F4").Value = "CMC"
$ws.Range("G4").Value = "CMC"

$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "23"
$ws.Range("H4").Style = "Normal"

$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "1111111"
$ws.Range("I4").Style = "Normal"
